# The "Cases" query stored in B2 (startup sheet) incorrectly returned an
# extra `Cohort` column that wasn't part of the intended result set.
# Remove the erroneous trailing `co.cohort_description` return line from
# the Cypher query text (fixing the query error mentioned in the commit
# message "Fixed variables and query errors ...").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newCasesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Dalmatian']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

# Strip the trailing newline left by the here-string so the cell text
# ends right after "Response to Treatment`" (no trailing blank line).
$newCasesQuery = $newCasesQuery.TrimEnd("`r", "`n")

$ws.Range("B2").Value = $newCasesQuery

# Reflect the resulting selection/active cell after the edit.
$ws.Activate()
$ws.Range("B2").Select()

$wb.Save()
